# Automatizzazione del file ottimizzazione
#
# Applies the authored edits to the "main" sheet of the workbook:
#   - C8  (Banana yield)              : 15   -> 1
#   - C9  (Price of Banana)           : 7    -> formula =162-7  (155)
#   - K9  reference text/link         : "Banana-coffee system cropping guide"
#                                        -> "numbeo" (new hyperlink source)
#   - C18 (Percentage of smallholders
#          to be covered)             : 1    -> 1.471149435534613E-4
#          (back-solved so that C29, Total required forestry investment,
#           comes out to ~0.9 MSh)
#
# All downstream formula cells (C20, C23-C35 on "main", plus the lookup
# formulas on sheets S / Y / Z / VA that reference main!C29/C30/C32-C35)
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Banana yield (kg/plant): 15 -> 1
$ws.Range("C8").Value = 1

# Price of Banana (Sh/kg): was a flat 7, now derived from a formula
$ws.Range("C9").Formula = "=162-7"

# Replace the reference note in K9 with the new "numbeo" source and
# point it at a fresh hyperlink (the old "Banana-coffee system cropping
# guide" reference text is reused elsewhere in the sheet, so this cell's
# own string is simply swapped for the new one).
$ws.Range("K9").Value = "numbeo"
$ws.Hyperlinks.Add($ws.Range("K9"), "https://www.numbeo.com/cost-of-living/country_result.jsp?country=Kenya")

# Percentage of the smallholders to be covered - solved so that the
# "Total required forestry investment" (C29) lands on 0.9 MSh.
$ws.Range("C18").Value = [double]"1.471149435534613E-4"

# Leave the selection on the cell that was just edited, matching the
# frozen-pane scroll position recorded after the edit.
$ws.Range("C9").Select()
